# The deck originally had 5 slides:
#   1. CSD 122 Final Project / RPN Logic Calculator (title slide)
#   2. RPN Logic Calculators (overview)
#   3. Example calculation 1   <- to be removed
#   4. Example calculation 2   <- to be removed
#   5. Q & A (closing slide)
#
# Per the commit message ("The powerpoint is supposed to be 3 slides, and I
# had 5 ... edited original project ppt down to 3 slides, and put the extra
# slides into a separate ppt"), slides 3 and 4 ("Example calculation 1" and
# "Example calculation 2") are removed, leaving the title slide, the
# overview slide and the closing Q&A slide (now slide 3) in place, with
# their content, shapes and relationships left untouched.

$p = $ppt.ActivePresentation

# Walk the deck and delete the two "Example calculation" slides by matching
# their title text, rather than hard-coding indexes, so the script is
# resilient to any re-ordering.
$titlesToRemove = @('Example calculation 1', 'Example calculation 2')

for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $title = ""
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $title = $shape.TextFrame.TextRange.Text
            break
        }
    }
    foreach ($t in $titlesToRemove) {
        if ($title -eq $t) {
            $slide.Delete()
            break
        }
    }
}

Write-Output "Final slide count: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    Write-Output "Slide $i : $($s.Shapes.Item(1).TextFrame.TextRange.Text)"
}
